# Applies the changes described by the commit diff:
#  - Update several cell values (and number formats) on the "Schedule" sheet
#  - Update the active selection on the "Schedule" sheet
#  - Update the active selection on the "Playoffs" sheet

$wb = $excel.ActiveWorkbook

$schedule = $wb.Worksheets.Item("Schedule")

# --- Simple value-only edits (format / style unchanged) ---
$schedule.Range("C4").Value  = 125.3
$schedule.Range("E7").Value  = 120.7
$schedule.Range("C9").Value  = 145
$schedule.Range("C14").Value = 115
$schedule.Range("C16").Value = 132.05454545454543
$schedule.Range("C18").Value = 121
$schedule.Range("C19").Value = 130
$schedule.Range("E19").Value = 123

# --- Edits where the number format also changes from 0.00 to 0.0 ---
$schedule.Range("E9").NumberFormat  = "0.0"
$schedule.Range("E9").Value         = 132.05454545454543

$schedule.Range("C10").NumberFormat = "0.0"
$schedule.Range("C10").Value        = 140

$schedule.Range("E10").NumberFormat = "0.0"
$schedule.Range("E10").Value        = 132.92909090909092

$schedule.Range("C11").NumberFormat = "0.0"
$schedule.Range("C11").Value        = 133

$schedule.Range("E11").NumberFormat = "0.0"
$schedule.Range("E11").Value        = 123.8

$schedule.Range("C12").NumberFormat = "0.0"
$schedule.Range("C12").Value        = 137.28181818181818

$schedule.Range("C13").NumberFormat = "0.0"
$schedule.Range("C13").Value        = 125.37

$schedule.Range("E13").NumberFormat = "0.0"
$schedule.Range("E13").Value        = 125.06

$schedule.Range("C17").NumberFormat = "0.0"
$schedule.Range("C17").Value        = 132.92909090909092

$schedule.Range("E18").NumberFormat = "0.0"
$schedule.Range("E18").Value        = 100

# --- Value-only edit, format already 0.0 ---
$schedule.Range("E12").Value = 132.5

# --- Update view / selection state ---
$schedule.Activate() | Out-Null
$schedule.Range("E18").Select() | Out-Null

$playoffs = $wb.Worksheets.Item("Playoffs")
$playoffs.Activate() | Out-Null
$playoffs.Range("E12").Select() | Out-Null

# Re-activate the Schedule sheet so it stays the visible tab, matching the
# workbook's activeTab setting.
$schedule.Activate() | Out-Null
